$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("G3").Value = "2016-08-20 18:57:17"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("H3").Value = "2016-08-20 18:57:13"
$ws2.Range("K3").Value = "2016-08-20 18:57:29"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("H3").Value = "2016-08-20 18:57:17"
$ws3.Range("K3").Value = "2016-08-20 18:57:35"
